$d = $word.ActiveDocument

# Locate the paragraph that holds the manual page break run (Chr(12)). It
# currently sits directly after the LASSO Model Params block and directly
# before the "LOGG REGG" section. We insert the new "LOG REG" params block
# right before that paragraph, so the manual page break ends up attached to
# the end of the new block (matching the target structure) instead of
# sitting in its own empty paragraph.
function Find-PageBreakParagraph($doc) {
    $cnt = $doc.Paragraphs.Count
    for ($i = 1; $i -le $cnt; $i++) {
        $para = $doc.Paragraphs.Item($i)
        $t = $para.Range.Text
        if ($t.Length -ge 1 -and [int][char]$t[0] -eq 12) {
            return $para
        }
    }
    return $null
}

$pageBreakPara = Find-PageBreakParagraph $d
if ($pageBreakPara -eq $null) {
    throw "Could not locate the page-break paragraph"
}

# All eight lines (one blank, the "LOG REG" heading, and six data lines) are
# inserted as brand-new paragraphs ahead of the page-break paragraph.
$lines = @(
    "",
    "LOG REG",
    "[1.04498522] [[-6.95709777e-04 -2.90027622e-01  2.08937195e-03  1.32557324e-01",
    "   4.06217316e-03  9.07083071e-01  2.00912133e-02  1.27812271e+00",
    "  -6.95709777e-04 -2.90027622e-01  3.99054651e-03 -3.81511254e-01",
    "   1.26657261e-03 -3.52195135e-01  5.97542420e-03 -5.98872523e-01",
    "   1.91568790e-03 -2.39266763e-02 -5.08093160e-03 -6.56741758e-01",
    "  -6.16172984e-03 -2.83947589e-01  2.55854632e-03  2.04468335e+00]]"
)

$insertText = ($lines -join "`r") + "`r"
$pageBreakPara.Range.InsertBefore($insertText)

# The paragraph collection shifted, so re-locate the page-break paragraph,
# then fold its immediately preceding paragraph (the last inserted data
# line) into it by deleting the paragraph mark that separates them. That
# leaves the manual page break as the tail of that last data-line paragraph.
$pageBreakPara = Find-PageBreakParagraph $d
$prevPara = $pageBreakPara.Previous()
$markRange = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
$markRange.Delete()

Write-Output "Done. New paragraph count: $($d.Paragraphs.Count)"
